$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  "B2" = "92.14%"
  "C2" = "92.14%"
  "D2" = "92.14%"
  "E2" = "92.15%"
  "F2" = "92.14%"
  "G2" = "92.14%"
  "H2" = "92.17%"
  "I2" = "91.98%"
  "J2" = "91.52%"
  "K2" = "65.59%"
  "B3" = "92.14%"
  "C3" = "92.14%"
  "D3" = "92.15%"
  "E3" = "92.15%"
  "F3" = "92.15%"
  "G3" = "92.13%"
  "H3" = "92.09%"
  "I3" = "92.02%"
  "J3" = "90.87%"
  "K3" = "81.72%"
  "B4" = "92.14%"
  "C4" = "92.15%"
  "D4" = "92.14%"
  "E4" = "92.15%"
  "F4" = "92.10%"
  "G4" = "92.04%"
  "H4" = "91.95%"
  "I4" = "89.89%"
  "J4" = "53.93%"
  "K4" = "50.26%"
  "B5" = "92.15%"
  "C5" = "92.14%"
  "D5" = "92.13%"
  "E5" = "92.14%"
  "F5" = "92.15%"
  "G5" = "92.15%"
  "H5" = "92.14%"
  "I5" = "92.09%"
  "J5" = "91.91%"
  "K5" = "84.93%"
  "B6" = "92.14%"
  "C6" = "92.14%"
  "D6" = "92.14%"
  "E6" = "92.17%"
  "F6" = "92.11%"
  "G6" = "92.16%"
  "H6" = "92.02%"
  "I6" = "91.38%"
  "J6" = "73.82%"
  "K6" = "50.12%"
  "B7" = "92.14%"
  "C7" = "92.14%"
  "D7" = "92.15%"
  "E7" = "92.14%"
  "F7" = "92.14%"
  "G7" = "92.15%"
  "H7" = "92.15%"
  "I7" = "92.15%"
  "J7" = "92.14%"
  "K7" = "92.14%"
  "B8" = "92.14%"
  "C8" = "92.15%"
  "D8" = "92.14%"
  "E8" = "92.14%"
  "F8" = "92.14%"
  "G8" = "92.15%"
  "H8" = "92.14%"
  "I8" = "92.14%"
  "J8" = "92.14%"
  "K8" = "92.15%"
}

foreach ($key in $data.Keys) {
  $ws.Range($key).Formula = "=""" + $data[$key] + """"
}

$rng = $ws.Range("B2:K8")
$rng.Copy()
$rng.PasteSpecial(-4163)
$ws.Range("A1").Select()
